$wb = $excel.ActiveWorkbook

# --- HPEbP sheet: correct the electrolysis energy-input formula -----------
# The "+46" term is dropped from the denominator (118/(162+2+46) -> 118/(162+2))
$wsH = $wb.Worksheets.Item("HPEbP")
$wsH.Range("B3").Formula = "=118/(162+2)"

# --- About sheet: tag the file with the state/region and an effective date -
$wsA = $wb.Worksheets.Item("About")
$wsA.Range("B1").Value = "Minnesota"
$wsA.Range("C1").Value = 44811
$wsA.Range("C1").NumberFormat = "mm-dd-yy"

# --- Restore each sheet's last-used selection / active view ---------------
$wsA.Activate()
$wsA.Range("B14").Select()

$wsI = $wb.Worksheets.Item("IEA Data")
$wsI.Activate()
$wsI.Range("D7:F7").Select()

$wsH.Activate()
$wsH.Range("C3").Select()
